$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: "Amazonaws@23" -> "Mayesha@1" and turn it into a hyperlink
$ws.Hyperlinks.Add($ws.Range("B2"), "https://mayesha1", "", "", "Mayesha@1")

# Update A3 text stays the same value "mosajgohar2" (no actual change needed,
# the shared-string reordering is incidental) - leave as-is.

# Move the active selection to F9 (matches the final sheetView selection)
$ws.Range("F9").Select()
